# Adds May 2021 data (rows 419-449, dates 44317-44347) to each of the five
# COVID tracking sheets of the "Reggio Calabria" workbook:
#   - column A: one row per day from 01/05/2021 (44317) to 31/05/2021 (44347)
#   - column C: the actual daily count is only known for the first three new
#     days (rows 419-421); later rows are placeholders with just the date
#   - column D: rolling 7-day AVERAGE(C(row-6):Crow), continuing the pattern
#     already used on rows 9-418
#   - column E (Ricoveri/sheet4 only): day-over-day delta C(row)-C(row-1)
# Also moves the active-tab/tabSelected pointer from "Dimessi   Guariti"
# (sheet 3) to "Terapia" (sheet 5), and updates each sheet's selection to
# C419:C421 (the newly entered numeric cells), matching what happens when a
# person types new values into those cells and leaves the selection there.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Per-sheet plan: 1-based worksheet index -> hashtable of C419/C420/C421
# values, and whether the sheet also carries an "E" (delta) column.
$plans = @(
    @{ Index = 1; C = @{ 419 = 49;  420 = 79;  421 = 120 }; HasE = $false },
    @{ Index = 2; C = @{ 419 = 1;   420 = 2;   421 = 2   }; HasE = $false },
    @{ Index = 3; C = @{ 419 = 120; 420 = 107; 421 = 122 }; HasE = $false },
    @{ Index = 4; C = @{ 419 = 129; 420 = 118; 421 = 125 }; HasE = $true  },
    @{ Index = 5; C = @{ 419 = 8;   420 = 9;   421 = 8   }; HasE = $false }
)

foreach ($plan in $plans) {
    $ws = $wb.Worksheets.Item($plan.Index)

    # Column A: consecutive dates for every new row, 419 (01/05/2021,
    # serial 44317) through 449 (31/05/2021, serial 44347).
    for ($row = 419; $row -le 449; $row++) {
        $ws.Cells.Item($row, 1).Value2 = 44317 + ($row - 419)
    }
    # Keep column A's date format identical to the existing rows above.
    $ws.Range("A418").Copy()
    $ws.Range("A419:A449").PasteSpecial($xlPasteFormats)

    # Column C + rolling 7-day average in D, only for the three rows that
    # actually carry a reported count (419-421).
    foreach ($row in 419..421) {
        $ws.Cells.Item($row, 3).Value2 = $plan.C[$row]
        $avgFirstRow = $row - 6
        $ws.Cells.Item($row, 4).Formula = "=AVERAGE(C$($avgFirstRow):C$row)"
        if ($plan.HasE) {
            $prevRow = $row - 1
            $ws.Cells.Item($row, 5).Formula = "=C$row-C$prevRow"
        }
    }
    # Match the number format/font of the averages (and deltas) column to
    # the existing cells directly above the new block.
    $ws.Range("D418").Copy()
    $ws.Range("D419:D421").PasteSpecial($xlPasteFormats)
    if ($plan.HasE) {
        $ws.Range("E418").Copy()
        $ws.Range("E419:E421").PasteSpecial($xlPasteFormats)
    }

    # Move the selection to the freshly entered values, as if the user had
    # just finished typing them in and left the cursor there.
    $ws.Range("C419:C421").Select()
}

# The active sheet moves from "Dimessi   Guariti" (index 3) to "Terapia"
# (index 5); selecting it last also updates workbook-level activeTab and
# moves tabSelected onto this sheet's view.
$wsTerapia = $wb.Worksheets.Item(5)
$wsTerapia.Select()
$wsTerapia.Range("C419:C421").Select()
